# "Colocando header nos gráficos"
# Adds an A1 header label to each data table (matching the style already
# used by B1:E1), drops the bold/border style from the A-column data rows
# (keeping it only on the new header row), fixes a handful of accented
# Portuguese labels, removes the now-unused "Teto" row from the
# Emissoes Totais sheet, and refreshes the Custo Total sheet's header/values.

$wb = $excel.ActiveWorkbook

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Sheets 1-4: "Potencia Acumulada", "Geracao Periodo Medio",
# "Atendimento a Ponta", "Potencia Incremental" all share the same
# row layout (Fonte/Tecnologia header + 11 source rows).
# ---------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # New header cell in A1, styled like the rest of row 1 (B1:E1).
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Fix accented labels.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # The header style moves up to row 1 only - strip it from A2:A12.
    $ws.Range("A2:A12").ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy() | Out-Null
$ws5.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Range("A2:A3").ClearFormats()

# Row 4 ("Teto") is no longer used.
$ws5.Rows("4:4").Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy() | Out-Null
$ws6.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$ws6.Range("A1").Value = "Tipo Expansão"

# B1 becomes a text "2015" (like the year headers on the other sheets),
# not a number - force text type with a quote prefix, then restore the
# normal header style (border/bold/centered) that the prefix step
# disturbed.
$ws6.Cells.Item(1, 2).Value = "'2015"
$ws6.Range("A2").Copy() | Out-Null
$ws6.Range("B1").PasteSpecial($xlPasteFormats) | Out-Null

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 571

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99

$ws6.Range("A2:A3").ClearFormats()

Write-Output "edit applied"
